$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "63-27=36"
$t.Cell(1,2).Range.Text = "58+26=84"
$t.Cell(1,3).Range.Text = "49+5=54"
$t.Cell(1,4).Range.Text = "38-29=9"
$t.Cell(1,5).Range.Text = "84-7=77"

$t.Cell(2,1).Range.Text = "82-19=63"
$t.Cell(2,2).Range.Text = "76-8=68"
$t.Cell(2,3).Range.Text = "46-7=39"
$t.Cell(2,4).Range.Text = "71-3=68"
$t.Cell(2,5).Range.Text = "11-3=8"

$t.Cell(3,1).Range.Text = "60-8=52"
$t.Cell(3,2).Range.Text = "50-2=48"
$t.Cell(3,3).Range.Text = "59+28=87"
$t.Cell(3,4).Range.Text = "90-28=62"
$t.Cell(3,5).Range.Text = "64-6=58"

$t.Cell(4,1).Range.Text = "68-19=49"
$t.Cell(4,2).Range.Text = "69+6=75"
$t.Cell(4,3).Range.Text = "26+56=82"
$t.Cell(4,4).Range.Text = "21-7=14"
$t.Cell(4,5).Range.Text = "94-76=18"

$t.Cell(5,1).Range.Text = "33+38=71"
$t.Cell(5,2).Range.Text = "41-28=13"
$t.Cell(5,3).Range.Text = "76-39=37"
$t.Cell(5,4).Range.Text = "76+6=82"
$t.Cell(5,5).Range.Text = "81-53=28"

$t.Cell(6,1).Range.Text = "60-5=55"
$t.Cell(6,2).Range.Text = "59+8=67"
$t.Cell(6,3).Range.Text = "24-7=17"
$t.Cell(6,4).Range.Text = "57-9=48"
$t.Cell(6,5).Range.Text = "53-46=7"

$t.Cell(7,1).Range.Text = "58+15=73"
$t.Cell(7,2).Range.Text = "19+47=66"
$t.Cell(7,3).Range.Text = "49+29=78"
$t.Cell(7,4).Range.Text = "93-78=15"
$t.Cell(7,5).Range.Text = "72-5=67"

$t.Cell(8,1).Range.Text = "37+47=84"
$t.Cell(8,2).Range.Text = "53+29=82"
$t.Cell(8,3).Range.Text = "40-17=23"
$t.Cell(8,4).Range.Text = "43-7=36"
$t.Cell(8,5).Range.Text = "58+36=94"

$t.Cell(9,1).Range.Text = "74-67=7"
$t.Cell(9,2).Range.Text = "73-29=44"
$t.Cell(9,3).Range.Text = "51-14=37"
$t.Cell(9,4).Range.Text = "5+69=74"
$t.Cell(9,5).Range.Text = "70-2=68"

$t.Cell(10,1).Range.Text = "26+49=75"
$t.Cell(10,2).Range.Text = "56+38=94"
$t.Cell(10,3).Range.Text = "96-68=28"
$t.Cell(10,4).Range.Text = "39+16=55"
$t.Cell(10,5).Range.Text = "8+88=96"

$t.Cell(11,1).Range.Text = "17+18=35"
$t.Cell(11,2).Range.Text = "25+26=51"
$t.Cell(11,3).Range.Text = "37+27=64"
$t.Cell(11,4).Range.Text = "48+48=96"
$t.Cell(11,5).Range.Text = "74-9=65"

$t.Cell(12,1).Range.Text = "17+46=63"
$t.Cell(12,2).Range.Text = "59+17=76"
$t.Cell(12,3).Range.Text = "3+48=51"
$t.Cell(12,4).Range.Text = "56+35=91"
$t.Cell(12,5).Range.Text = "7+29=36"

$t.Cell(13,1).Range.Text = "18+54=72"
$t.Cell(13,2).Range.Text = "29+34=63"
$t.Cell(13,3).Range.Text = "65+6=71"
$t.Cell(13,4).Range.Text = "5+36=41"
$t.Cell(13,5).Range.Text = "54-35=19"

$t.Cell(14,1).Range.Text = "22-4=18"
$t.Cell(14,2).Range.Text = "76+9=85"
$t.Cell(14,3).Range.Text = "91-13=78"
$t.Cell(14,4).Range.Text = "76-7=69"
$t.Cell(14,5).Range.Text = "58+19=77"

$t.Cell(15,1).Range.Text = "63-45=18"
$t.Cell(15,2).Range.Text = "4+18=22"
$t.Cell(15,3).Range.Text = "5+37=42"
$t.Cell(15,4).Range.Text = "15+17=32"
$t.Cell(15,5).Range.Text = "66-17=49"

$t.Cell(16,1).Range.Text = "32-26=6"
$t.Cell(16,2).Range.Text = "65+6=71"
$t.Cell(16,3).Range.Text = "82-64=18"
$t.Cell(16,4).Range.Text = "94-59=35"
$t.Cell(16,5).Range.Text = "86-7=79"

$t.Cell(17,1).Range.Text = "69+25=94"
$t.Cell(17,2).Range.Text = "46+28=74"
$t.Cell(17,3).Range.Text = "26+47=73"
$t.Cell(17,4).Range.Text = "29+26=55"
$t.Cell(17,5).Range.Text = "71-38=33"

$t.Cell(18,1).Range.Text = "97-18=79"
$t.Cell(18,2).Range.Text = "91-43=48"
$t.Cell(18,3).Range.Text = "35+49=84"
$t.Cell(18,4).Range.Text = "50-29=21"
$t.Cell(18,5).Range.Text = "50-8=42"

$t.Cell(19,1).Range.Text = "31-13=18"
$t.Cell(19,2).Range.Text = "48+18=66"
$t.Cell(19,3).Range.Text = "8+39=47"
$t.Cell(19,4).Range.Text = "35+19=54"
$t.Cell(19,5).Range.Text = "86-18=68"

$t.Cell(20,1).Range.Text = "87-39=48"
$t.Cell(20,2).Range.Text = "88-69=19"
$t.Cell(20,3).Range.Text = "15+6=21"
$t.Cell(20,4).Range.Text = "62-58=4"
$t.Cell(20,5).Range.Text = "38+9=47"
